# Updates the cryptos price list (Price and Volume(1h) columns) with the
# latest scraped figures. Two pairs of rows also had their coin data
# swapped back into rank order (EnergySwap/NEARProtocol at rows 35-36,
# and Cosmos/ONDO at rows 48-49) -- row numbers/ranks (column A) stay put,
# only the Coin/Link/Price/Volume cells move.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values that look like plain decimal numbers (e.g. "1.00", "0.430")
# must be written to cells pre-formatted as Text, otherwise Excel's
# automatic type detection would convert them to numbers and silently
# drop significant trailing/leading zeros (e.g. "1.00" -> 1, "0.430" -> 0.43).

$ws.Range("D2").Value = '56.405.71'
$ws.Range("E2").Value = '  -0.10%  '

$ws.Range("D3").Value = '3.013.57'
$ws.Range("E3").Value = '  +2.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '509.52'
$ws.Range("E5").Value = '  +2.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.04'
$ws.Range("E6").Value = '  +4.11%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.430'
$ws.Range("E8").Value = '  +1.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.10'
$ws.Range("E9").Value = '  -0.87%  '

$ws.Range("E10").Value = '  +1.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.368'
$ws.Range("E11").Value = '  +4.90%  '

$ws.Range("D12").Value = '3.527.02'
$ws.Range("E12").Value = '  +1.94%  '

$ws.Range("E13").Value = '  +0.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.34'
$ws.Range("E14").Value = '  -2.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000162'
$ws.Range("E15").Value = '  +3.39%  '

$ws.Range("D16").Value = '56.433.57'
$ws.Range("E16").Value = '  -0.16%  '

$ws.Range("D17").Value = '3.015.06'
$ws.Range("E17").Value = '  +2.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.92'
$ws.Range("E18").Value = '  -1.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.07'
$ws.Range("E19").Value = '  +4.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.00'
$ws.Range("E20").Value = '  +3.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '331.83'
$ws.Range("E21").Value = '  +4.68%  '

$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.499'
$ws.Range("E23").Value = '  +3.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.65'
$ws.Range("E24").Value = '  +3.71%  '

$ws.Range("D25").Value = '3.137.86'
$ws.Range("E25").Value = '  +2.15%  '

$ws.Range("E26").Value = '  -0.38%  '

$ws.Range("E27").Value = '  +2.62%  '

$ws.Range("D28").Value = '0.0₃0939'
$ws.Range("E28").Value = '  +8.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.36'
$ws.Range("E29").Value = '  -1.85%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.80'
$ws.Range("E30").Value = '  -2.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.79'
$ws.Range("E31").Value = '  +2.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.17'
$ws.Range("E32").Value = '  +3.10%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.36'
$ws.Range("E33").Value = '  +2.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '152.78'
$ws.Range("E34").Value = '  +0.07%  '

$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.46'
$ws.Range("E35").Value = '  -0.13%  '

$ws.Range("B36").Value = 'EnergySwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '27.02'
$ws.Range("E36").Value = '  +14.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.81'
$ws.Range("E37").Value = '  +2.52%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.22'
$ws.Range("E38").Value = '  +1.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0662'
$ws.Range("E39").Value = '  +1.52%  '

$ws.Range("D40").Value = '3.056.13'
$ws.Range("E40").Value = '  +2.68%  '

$ws.Range("E41").Value = '  -2.56%  '

$ws.Range("E42").Value = '  +0.20%  '

$ws.Range("E43").Value = '  +3.27%  '

$ws.Range("E44").Value = '  +2.63%  '

$ws.Range("D45").Value = '2.203.94'
$ws.Range("E45").Value = '  +2.85%  '

$ws.Range("E46").Value = '  -0.41%  '

$ws.Range("E47").Value = '  +4.99%  '

$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.924'
$ws.Range("E48").Value = '  +0.74%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.83'
$ws.Range("E49").Value = '  -0.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.57'
$ws.Range("E50").Value = '  +3.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0856'
$ws.Range("E51").Value = '  -0.01%  '
